$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9 (shifts old rows 9-10 down to 10-11) to
# hold the new "Request time off work due to domestic violence" entry in
# its alphabetically-sorted position.
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "Request time off work due to domestic violence"
$ws.Range("B9").Value = "https://www.illinoislegalaid.org/legal-information/request-time-work-due-domestic-abuse-letter"

# Rebuild every hyperlink (the row insert does not shift the existing
# hyperlink anchors automatically) so each link ends up pointing at its
# correct, possibly-shifted, cell.
$ws.Range("A1:B11").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.illinoislegalaid.org/legal-information/appearance")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://www.illinoislegalaid.org/legal-information/fee-waiver")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.illinoislegalaid.org/legal-information/collection-proof-debtor-letter")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://www.illinoislegalaid.org/legal-information/request-collection-agency-stop-contacting")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://www.illinoislegalaid.org/legal-information/end-illegal-lockout-demand")
$ws.Hyperlinks.Add($ws.Range("B10"), "https://www.illinoislegalaid.org/legal-information/security-deposit-demand-letter")
$ws.Hyperlinks.Add($ws.Range("B8"), "https://www.illinoislegalaid.org/legal-information/housing-discrimination-complaint-idhr")
$ws.Hyperlinks.Add($ws.Range("B11"), "https://www.illinoislegalaid.org/legal-information/stop-wage-assignment-letter")
$ws.Hyperlinks.Add($ws.Range("B9"), "https://www.illinoislegalaid.org/legal-information/request-time-work-due-domestic-abuse-letter")

# Restore the hyperlink cell style (Hyperlinks.Add bumps in a fresh style
# index; put every linked cell back on the shared "Hyperlink" style).
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B7").Style = "Hyperlink"
$ws.Range("B4").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("B5").Style = "Hyperlink"
$ws.Range("B10").Style = "Hyperlink"
$ws.Range("B8").Style = "Hyperlink"
$ws.Range("B11").Style = "Hyperlink"
$ws.Range("B9").Style = "Hyperlink"

# Re-apply the A2:A10 name sort so the worksheet's remembered sort range
# grows to include the new row.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A10"))
$ws.Sort.SetRange($ws.Range("A2:B10"))
$ws.Sort.Apply()

# Column A widened (best-fit) to accommodate the longest new label.
$ws.Columns.Item(1).ColumnWidth = 43.8

# Match the author's last selected cell.
[void]$ws.Range("B20").Select()
